# Replace the lattice-multiplication exercise cells with a new set of
# problems. Each table cell holds 5 lines joined by manual line breaks
# (w:br, represented as Chr(11) in Word's Range.Text):
#   "AB x CD" / "  C    D" / "  ----" / "A|    |" / "B|    |"
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "12 x 16" + [char]11 + "  1    6" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "2|    |"
$t.Cell(1,2).Range.Text = "52 x 13" + [char]11 + "  1    3" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "2|    |"
$t.Cell(1,3).Range.Text = "99 x 96" + [char]11 + "  9    6" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "9|    |"
$t.Cell(2,1).Range.Text = "63 x 41" + [char]11 + "  4    1" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "3|    |"
$t.Cell(2,2).Range.Text = "29 x 34" + [char]11 + "  3    4" + [char]11 + "  ----" + [char]11 + "2|    |" + [char]11 + "9|    |"
$t.Cell(2,3).Range.Text = "53 x 62" + [char]11 + "  6    2" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "3|    |"
$t.Cell(3,1).Range.Text = "87 x 96" + [char]11 + "  9    6" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "7|    |"
$t.Cell(3,2).Range.Text = "66 x 57" + [char]11 + "  5    7" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "6|    |"
$t.Cell(3,3).Range.Text = "16 x 82" + [char]11 + "  8    2" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "6|    |"
$t.Cell(4,1).Range.Text = "31 x 30" + [char]11 + "  3    0" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "1|    |"
$t.Cell(4,2).Range.Text = "95 x 21" + [char]11 + "  2    1" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "5|    |"
$t.Cell(4,3).Range.Text = "72 x 38" + [char]11 + "  3    8" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "2|    |"
$t.Cell(5,1).Range.Text = "18 x 29" + [char]11 + "  2    9" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "8|    |"
$t.Cell(5,2).Range.Text = "41 x 20" + [char]11 + "  2    0" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "1|    |"
$t.Cell(5,3).Range.Text = "62 x 50" + [char]11 + "  5    0" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "2|    |"
